$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update parameter values for rows 2-13
$ws.Range("B2").Value = -0.06807742355664341
$ws.Range("B3").Value = -0.06858330716945422
$ws.Range("B4").Value = -0.1048294541821364
$ws.Range("B5").Value = -0.06851356552895799
$ws.Range("B6").Value = -0.2983606849676014
$ws.Range("B7").Value = 0.137610500119417
$ws.Range("B8").Value = 0.04065523770289773
$ws.Range("B9").Value = 0.0007620119862800202
$ws.Range("B10").Value = -0.04983344429657881
$ws.Range("B11").Value = -0.1147592499693111
$ws.Range("B12").Value = -0.217234202289305
$ws.Range("B13").Value = -0.01906499818396121

# Remove the last row (shot_during_regular_play), shifting cells up
$ws.Range("A14:B14").Delete()
